# Update the LR-pairs data table with new TPM-based values.
# Rows 2-16 get updated values for existing Sending/Target cluster combos,
# and rows 17-21 are newly added rows for the "Resolving-Mac" sending cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Il7", "Il2rg", "ECs", 2, 0.6666666666666666, 0.4161346666666667, 1.248404, 0.4872818821417071, 0.487281882141707, 3, 1, 50.26671733333333, 150.800152, 0.2619764206727233, 0.2619764206727234, 20.91772366193423, 188.259512957408, 0.1276563633421522, 0.1276563633421522)
    ,@("ECs", "Il7", "Il2rg", "FAPs", 2, 0.6666666666666666, 0.4161346666666667, 1.248404, 0.4872818821417071, 0.487281882141707, 3, 1, 3.113353333333333, 9.34006, 0.01622594841727, 0.01622594841727001, 1.295574251582222, 11.66016826424, 0.00790661068430158, 0.007906610684301582)
    ,@("ECs", "Il7", "Il2rg", "Inflammatory-Mac", 2, 0.6666666666666666, 0.4161346666666667, 1.248404, 0.4872818821417071, 0.487281882141707, 3, 1, 80.63290666666667, 241.89872, 0.4202367172077739, 0.4202367172077739, 33.55414773809778, 301.98732964288, 0.2047737385060564, 0.2047737385060563)
    ,@("ECs", "Il7", "Il2rg", "MuSCs", 2, 0.6666666666666666, 0.4161346666666667, 1.248404, 0.4872818821417071, 0.487281882141707, 3, 1, 2.274154, 6.822462, 0.01185227038057408, 0.01185227038057408, 0.9463543167386667, 8.517188850648001, 0.005775396618698545, 0.005775396618698545)
    ,@("ECs", "Il7", "Il2rg", "Resolving-Mac", 2, 0.6666666666666666, 0.4161346666666667, 1.248404, 0.4872818821417071, 0.487281882141707, 3, 1, 55.58783666666667, 166.76351, 0.2897086433216586, 0.2897086433216586, 23.13202588200445, 208.18823293804, 0.1411697729904983, 0.1411697729904983)
    ,@("FAPs", "Il7", "Il2rg", "ECs", 3, 1, 0.298935, 0.8968050000000001, 0.3500443993403526, 0.3500443993403526, 3, 1, 50.26671733333333, 150.800152, 0.2619764206727233, 0.2619764206727234, 15.02648114604, 135.23833031436, 0.09170337881571897, 0.091703378815719)
    ,@("FAPs", "Il7", "Il2rg", "FAPs", 3, 1, 0.298935, 0.8968050000000001, 0.3500443993403526, 0.3500443993403526, 3, 1, 3.113353333333333, 9.34006, 0.01622594841727, 0.01622594841727001, 0.9306902786999999, 8.3762125083, 0.005679802367450824, 0.005679802367450825)
    ,@("FAPs", "Il7", "Il2rg", "Inflammatory-Mac", 3, 1, 0.298935, 0.8968050000000001, 0.3500443993403526, 0.3500443993403526, 3, 1, 80.63290666666667, 241.89872, 0.4202367172077739, 0.4202367172077739, 24.1039979544, 216.9359815896, 0.1471015092557569, 0.1471015092557569)
    ,@("FAPs", "Il7", "Il2rg", "MuSCs", 3, 1, 0.298935, 0.8968050000000001, 0.3500443993403526, 0.3500443993403526, 3, 1, 2.274154, 6.822462, 0.01185227038057408, 0.01185227038057408, 0.6798242259899999, 6.11841803391, 0.004148820866187507, 0.004148820866187507)
    ,@("FAPs", "Il7", "Il2rg", "Resolving-Mac", 3, 1, 0.298935, 0.8968050000000001, 0.3500443993403526, 0.3500443993403526, 3, 1, 55.58783666666667, 166.76351, 0.2897086433216586, 0.2897086433216586, 16.61714995395, 149.55434958555, 0.1014108880352385, 0.1014108880352385)
    ,@("MuSCs", "Il7", "Il2rg", "ECs", 2, 0.6666666666666666, 0.06745866666666667, 0.202376, 0.07899218376447857, 0.07899218376447857, 3, 1, 50.26671733333333, 150.800152, 0.2619764206727233, 0.2619764206727234, 3.390925729016889, 30.518331561152, 0.0206940895637401, 0.02069408956374011)
    ,@("MuSCs", "Il7", "Il2rg", "FAPs", 2, 0.6666666666666666, 0.06745866666666667, 0.202376, 0.07899218376447857, 0.07899218376447857, 3, 1, 3.113353333333333, 9.34006, 0.01622594841727, 0.01622594841727001, 0.2100226647288889, 1.89020398256, 0.001281723099129942, 0.001281723099129943)
    ,@("MuSCs", "Il7", "Il2rg", "Inflammatory-Mac", 2, 0.6666666666666666, 0.06745866666666667, 0.202376, 0.07899218376447857, 0.07899218376447857, 3, 1, 80.63290666666667, 241.89872, 0.4202367172077739, 0.4202367172077739, 5.439388373191111, 48.95449535872, 0.03319541599025769, 0.03319541599025769)
    ,@("MuSCs", "Il7", "Il2rg", "MuSCs", 2, 0.6666666666666666, 0.06745866666666667, 0.202376, 0.07899218376447857, 0.07899218376447857, 3, 1, 2.274154, 6.822462, 0.01185227038057408, 0.01185227038057408, 0.1534113966346667, 1.380702569712, 0.0009362367199285942, 0.0009362367199285943)
    ,@("MuSCs", "Il7", "Il2rg", "Resolving-Mac", 2, 0.6666666666666666, 0.06745866666666667, 0.202376, 0.07899218376447857, 0.07899218376447857, 3, 1, 55.58783666666667, 166.76351, 0.2897086433216586, 0.2897086433216586, 3.749881344417778, 33.74893209976, 0.02288471839142224, 0.02288471839142224)
    ,@("Resolving-Mac", "Il7", "Il2rg", "ECs", 2, 0.6666666666666666, 0.07146333333333334, 0.21439, 0.08368153475346168, 0.08368153475346167, 3, 1, 50.26671733333333, 150.800152, 0.2619764206727233, 0.2619764206727234, 3.592227176364445, 32.33004458728, 0.02192258895111199, 0.02192258895111199)
    ,@("Resolving-Mac", "Il7", "Il2rg", "FAPs", 2, 0.6666666666666666, 0.07146333333333334, 0.21439, 0.08368153475346168, 0.08368153475346167, 3, 1, 3.113353333333333, 9.34006, 0.01622594841727, 0.01622594841727001, 0.2224906070444444, 2.0024154634, 0.001357812266387656, 0.001357812266387656)
    ,@("Resolving-Mac", "Il7", "Il2rg", "Inflammatory-Mac", 2, 0.6666666666666666, 0.07146333333333334, 0.21439, 0.08368153475346168, 0.08368153475346167, 3, 1, 80.63290666666667, 241.89872, 0.4202367172077739, 0.4202367172077739, 5.762296286755556, 51.8606665808, 0.03516605345570298, 0.03516605345570298)
    ,@("Resolving-Mac", "Il7", "Il2rg", "MuSCs", 2, 0.6666666666666666, 0.07146333333333334, 0.21439, 0.08368153475346168, 0.08368153475346167, 3, 1, 2.274154, 6.822462, 0.01185227038057408, 0.01185227038057408, 0.1625186253533333, 1.46266762818, 0.0009918161757594345, 0.0009918161757594345)
    ,@("Resolving-Mac", "Il7", "Il2rg", "Resolving-Mac", 2, 0.6666666666666666, 0.07146333333333334, 0.21439, 0.08368153475346168, 0.08368153475346167, 3, 1, 55.58783666666667, 166.76351, 0.2897086433216586, 0.2897086433216586, 3.972492100988889, 35.7524289089, 0.02424326390449961, 0.02424326390449961)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$j]
    }
}
